$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A32").Value = "Matteo Zanlucchi"
$ws.Range("B32").Value = "Stefano Tita | Clitoriders"
$ws.Range("C32").Value = "Federico  Manica | iMontagna"
$ws.Range("D32").Value = "Luca Frasca | Clitoriders"
$ws.Range("E32").Value = "Federico  Fasanelli  | Herta Vernello"
$ws.Range("F32").Value = "Simone Miorelli | SBARX"
